$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
#          M (Precio promedio ponderado), P (Precio $/Kg)
# Values are permuted among data rows 2,3,4,5,7,8,9 (row 6 unchanged).
# New row values (captured from the target state):
$data = @{
    2 = @(44406, 160, 17000, 18000, 17500, 1346)
    3 = @(44379, 120, 12000, 13000, 12667, 974)
    4 = @(44469, 140, 13000, 14000, 13500, 1038)
    5 = @(44229, 120, 44000, 45000, 44500, 3423)
    7 = @(44159, 100, 23000, 24000, 23500, 1808)
    8 = @(44397, 140, 12500, 13000, 12750, 981)
    9 = @(44389, 120, 12000, 13000, 12500, 962)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
    $ws.Range("K$row").Value = $vals[2]
    $ws.Range("L$row").Value = $vals[3]
    $ws.Range("M$row").Value = $vals[4]
    $ws.Range("P$row").Value = $vals[5]
}
